$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The sheet used to hold a 3-row "Name_1/Name_2/Name_3" table in B2:D4
# (columns: Name | Age | IsVip). Replace it with a "Key/Value" dictionary
# rendering of a single data item (Name="Name_1", Value=25.7, IsVip=True),
# rendered twice side-by-side (columns B:C and E:F) for the two data-source
# panel providers under test.
$ws.Range("A1:D4").ClearContents()

# Panel 1 - columns B (Key) / C (Value)
$ws.Range("B1").Value = "Key"
$ws.Range("C1").Value = "Value"
$ws.Range("B2").Value = "Name"
$ws.Range("C2").Value = "Name_1"
$ws.Range("B3").Value = "Value"
$ws.Range("C3").Value = 25.7
$ws.Range("B4").Value = "IsVip"
$ws.Range("C4").Value = $true

# Panel 2 - columns E (Key) / F (Value) - identical content, second provider
$ws.Range("E1").Value = "Key"
$ws.Range("F1").Value = "Value"
$ws.Range("E2").Value = "Name"
$ws.Range("F2").Value = "Name_1"
$ws.Range("E3").Value = "Value"
$ws.Range("F3").Value = 25.7
$ws.Range("E4").Value = "IsVip"
$ws.Range("F4").Value = $true
